# Add 7 new Google-Forms style leaderboard responses as rows 75-81 on Sheet1.
# (Matches the target diff: new survey submissions appended below the
# existing responses, reusing the same row layout/formatting as the rows
# above.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New response rows: Timestamp (Excel serial), name, roll, profile URL, branch, pre_points
$newRows = @(
    @(45384.587546296294, "Piyush Dwivedi",            "B23354", "https://www.beecrowd.com.br/judge/en/profile/948818", "GE",  0),
    @(45384.592106481483, "Akshit Bhola",               "B23112", "https://www.beecrowd.com.br/judge/en/profile/949151", "CSE", 0),
    @(45384.592187499999, "Adit Raj",                   "B23480", "https://www.beecrowd.com.br/judge/en/profile/949154", "VLSI",0),
    @(45384.593553240738, "Alok Kumar Yadav",           "B23426", "https://www.beecrowd.com.br/judge/en/profile/949149", "ME",  0),
    @(45384.618969907409, "Rudraksh Rajendra Lande",    "B23176", "https://www.beecrowd.com.br/judge/en/profile/949172", "CSE", 0),
    @(45384.623391203706, "Yashodeep",                  "B23040", "https://www.beecrowd.com.br/judge/en/profile/948216", "MnC", 0),
    @(45384.633993055555, "Mohit Kumar",                "B23273", "https://www.beecrowd.com.br/judge/en/profile/948171", "EE",  0)
)

$startRow = 75
$templateRow = 73

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    # Bring over the same cell styling (borders / wrap / number format) that
    # every other response row already uses.
    $ws.Range("A" + $templateRow + ":F" + $templateRow).Copy()
    $ws.Range("A" + $r + ":F" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]

    $ws.Rows.Item($r).RowHeight = 27.6
}

$endRow = $startRow + $newRows.Count - 1

# Match the saved view state: scrolled down to the new rows with F75:F81 selected.
$ws.Range("F" + $startRow + ":F" + $endRow).Select()
